# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# Update the computed K values for rows 2-16 to reflect the recalculated figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 2
    4  = 0
    5  = 1
    6  = 0
    7  = 0
    8  = 2
    9  = 3
    10 = 3
    11 = 3
    12 = 1
    13 = 2
    14 = 7
    15 = 3
    16 = 4
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
